# Workbook currently has two sheets: "ODI Batting" (sheetId 1) and
# "ODI Bowling" (sheetId 2). This edit:
#   1. Inserts a brand-new "Player Info" sheet before "ODI Batting"
#      (becomes the first sheet; the others shift right).
#   2. Renames the "MATCH_CARD_LINK" column to "MATCH_CODE" on both the
#      "ODI Batting" and "ODI Bowling" sheets and rewrites the URL values
#      there down to the bare numeric match code.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "Player Info" sheet, inserted before "ODI Batting".
# ---------------------------------------------------------------------
$wsBattingForInsert = $wb.Worksheets.Item("ODI Batting")
$wsInfo = $wb.Worksheets.Add($wsBattingForInsert)
$wsInfo.Name = "Player Info"

# Header row formatting to match the other sheets' header style
# (bold, thin box border, centered / top aligned).
$infoHeader = $wsInfo.Range("A1:D1")
$infoHeader.Font.Bold = $true
$infoHeader.Borders.LineStyle = 1
$infoHeader.Borders.Weight = 2
$infoHeader.HorizontalAlignment = -4108
$infoHeader.VerticalAlignment = -4160

$wsInfo.Range("A1").Value = "ID"
$wsInfo.Range("B1").Value = "NAME"
$wsInfo.Range("C1").Value = "BATTING_HAND"
$wsInfo.Range("D1").Value = "BOWL_STYLE"

# Data row - ID is numeric-looking text, force it to stay text like the
# rest of this workbook's cells (they're all inlineStr / shared strings).
$wsInfo.Range("A2").NumberFormat = "@"
$wsInfo.Range("A2").Value = "5930"
$wsInfo.Range("A2").ClearFormats()

$wsInfo.Range("B2").Value = "Jacob Andrew Duffy"
$wsInfo.Range("C2").Value = "Right Handed"
$wsInfo.Range("D2").Value = "Right Arm Fast Medium"

# ---------------------------------------------------------------------
# Re-fetch the other two sheets by name now that the insert has shifted
# their positions - worksheet variables captured before "Add" resolve by
# position, not stable identity, so they'd otherwise point at the wrong
# sheet after the insert.
# ---------------------------------------------------------------------
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBowling = $wb.Worksheets.Item("ODI Bowling")

# ---------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (column D), URL -> code.
# ---------------------------------------------------------------------
$wsBatting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{ 2 = "4608"; 3 = "4625"; 4 = "4697" }
foreach ($row in $battingCodes.Keys) {
    $cell = $wsBatting.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $battingCodes[$row]
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE (column B), URL -> code.
# ---------------------------------------------------------------------
$wsBowling.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @{ 2 = "4608"; 3 = "4625"; 4 = "4697" }
foreach ($row in $bowlingCodes.Keys) {
    $cell = $wsBowling.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $bowlingCodes[$row]
    $cell.ClearFormats()
}
